$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 2 (Target cluster = "ECs"); this shifts the old row 3
# (Target cluster = "MuSCs") up into row 2, and the now-unused "ECs" shared
# string is dropped from the table automatically.
$ws.Rows(2).Delete()

# Update the recomputed TPM-derived metrics on the resulting row 2.
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.010147
$ws.Range("N2").Value = 3.030441
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.3513954595863334
$ws.Range("R2").Value = 3.162559136277001
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
